# =====================================================================
# Daily Report update: 2026-02-19 commit
#   - Appends a new date's worth of rows (date serial 46071) to
#     Daily_Data, one Registered/Eligible pair per depository.
#   - Refreshes the "latest day" rollups on Today_Summary (per
#     depository Eligible/Registered/Total_Stock).
#   - Refreshes the current-month (2026-02) rollups on Monthly_Stats,
#     both the top summary row and the per-depository detail rows.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Daily_Data: append rows 684-705 for date 46071
#    columns: Date, Region_Type, PREV_TOTAL, RECEIVED, WITHDRAWN,
#             NET_CHANGE, ADJUSTMENT, TOTAL_TODAY
# ---------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(684, "ASAHI DEPOSITORY LLC Registered", 0, 0, 0, 0, 0, 0),
    @(685, "ASAHI DEPOSITORY LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(686, "BRINK'S, INC. Registered", 73354.783, 0, 0, 0, -2079.184, 71275.599),
    @(687, "BRINK'S, INC. Eligible", 84460.738, 0, 0, 0, 2079.184, 86539.92200000001),
    @(688, "CNT DEPOSITORY, INC. Registered", 1246.06, 0, 0, 0, 0, 1246.06),
    @(689, "CNT DEPOSITORY, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(690, "DELAWARE DEPOSITORY Registered", 1633.941, 0, 0, 0, 0, 1633.941),
    @(691, "DELAWARE DEPOSITORY Eligible", 18459.584, 0, 0, 0, 0, 18459.584),
    @(692, "HSBC BANK, USA Registered", 1394.758, 0, 0, 0, 0, 1394.758),
    @(693, "HSBC BANK, USA Eligible", 9281.978999999999, 0, 0, 0, 0, 9281.978999999999),
    @(694, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 2395.448, 0, 0, 0, 0, 2395.448),
    @(695, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 0, 0, 0, 0, 0, 0),
    @(696, "JP MORGAN CHASE BANK NA Registered", 114061.421, 0, 0, 0, -459.633, 113601.788),
    @(697, "JP MORGAN CHASE BANK NA Eligible", 76408.66899999999, 0, 924.158, -924.158, 459.633, 75944.144),
    @(698, "LOOMIS INTERNATIONAL (US) LLC Registered", 61157.444, 0, 0, 0, -1947.656, 59209.788),
    @(699, "LOOMIS INTERNATIONAL (US) LLC Eligible", 69005.64, 0, 0, 0, 1947.656, 70953.296),
    @(700, "MALCA-AMIT USA, LLC Registered", 395.145, 0, 0, 0, 0, 395.145),
    @(701, "MALCA-AMIT USA, LLC Eligible", 0, 0, 0, 0, 0, 0),
    @(702, "MANFRA, TORDELLA & BROOKES, LLC Registered", 49920.248, 0, 0, 0, -1627.601, 48292.647),
    @(703, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 1804.683, 0, 0, 0, 1627.601, 3432.284),
    @(704, "STONEX PRECIOUS METALS LLC Registered", 14122.765, 0, 0, 0, 0, 14122.765),
    @(705, "STONEX PRECIOUS METALS LLC Eligible", 16.075, 0, 0, 0, 0, 16.075)
)

foreach ($row in $newRows) {
    $r = $row[0]

    $wsDaily.Cells.Item($r, 1).Value = 46071
    $wsDaily.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsDaily.Cells.Item($r, 2).Value = $row[1]
    $wsDaily.Cells.Item($r, 3).Value = $row[2]
    $wsDaily.Cells.Item($r, 4).Value = $row[3]
    $wsDaily.Cells.Item($r, 5).Value = $row[4]
    $wsDaily.Cells.Item($r, 6).Value = $row[5]
    $wsDaily.Cells.Item($r, 7).Value = $row[6]
    $wsDaily.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------
# 2) Today_Summary: refresh Eligible(B) / Registered(C) / Total_Stock(D)
#    for the depositories whose latest-day totals moved.
# ---------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("Today_Summary")

$todaySummaryUpdates = @(
    @(3, 86539.92200000001, 71275.599, $null),
    @(8, 75944.144, 113601.788, 189545.932),
    @(9, 70953.296, 59209.788, $null),
    @(11, 3432.284, 48292.647, $null)
)

foreach ($row in $todaySummaryUpdates) {
    $r = $row[0]
    $wsToday.Cells.Item($r, 2).Value = $row[1]
    $wsToday.Cells.Item($r, 3).Value = $row[2]
    if ($row[3] -ne $null) {
        $wsToday.Cells.Item($r, 4).Value = $row[3]
    }
}

# ---------------------------------------------------------------
# 3) Monthly_Stats: refresh the 2026-02 summary row and the
#    2026-02 per-depository detail rows (WITHDRAWN / TOTAL_TODAY).
# ---------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly_Stats")

# Top summary row for 2026-02 (Eligible / Registered / Grand_Total)
$wsMonthly.Cells.Item(2, 2).Value = 264627.284
$wsMonthly.Cells.Item(2, 3).Value = 313567.939
$wsMonthly.Cells.Item(2, 4).Value = 578195.223

# Per-depository detail rows (columns: YearMonth, Region_Type, RECEIVED,
# WITHDRAWN, TOTAL_TODAY). Only D (WITHDRAWN) and E (TOTAL_TODAY) move.
$monthlyDetailUpdates = @(
    @(10, $null, 86539.92200000001),
    @(11, $null, 71275.599),
    @(20, 924.158, 75944.144),
    @(21, $null, 113601.788),
    @(22, $null, 70953.296),
    @(23, $null, 59209.788),
    @(26, $null, 3432.284),
    @(27, $null, 48292.647)
)

foreach ($row in $monthlyDetailUpdates) {
    $r = $row[0]
    if ($row[1] -ne $null) {
        $wsMonthly.Cells.Item($r, 4).Value = $row[1]
    }
    $wsMonthly.Cells.Item($r, 5).Value = $row[2]
}
